$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("B08CS8YFK5")

$ws.Cells.Item(1, 1).Value = 'knee pads 661'
$ws.Cells.Item(2, 1).Value = 'knee pads ocr'
$ws.Cells.Item(3, 1).Value = 'knee pad inserts for tactical pants'
$ws.Cells.Item(4, 1).Value = 'knee pads jbm'
$ws.Cells.Item(5, 1).Value = 'knee protectors for toddlers'
$ws.Cells.Item(6, 1).Value = 'protec knee pads'
$ws.Cells.Item(7, 1).Value = 'nike thermal compression pants for men'
$ws.Cells.Item(8, 1).Value = 'under armour compression pants youth'
$ws.Cells.Item(9, 1).Value = 'under armour compression tights for men'
$ws.Cells.Item(10, 1).Value = 'capri pants adidas'
$ws.Cells.Item(11, 1).Value = 'capri pants exercise'
$ws.Cells.Item(12, 1).Value = 'capri pants nike'
$ws.Cells.Item(13, 1).Value = 'imucci knee pads'
$ws.Cells.Item(14, 1).Value = 'bb knee pads'
$ws.Cells.Item(15, 1).Value = 'knee pad and helmet'
$ws.Cells.Item(16, 1).Value = 'knee pad basketball kids'
$ws.Cells.Item(17, 1).Value = 'knee pad buttons'
$ws.Cells.Item(18, 1).Value = 'knee pad climbing'
$ws.Cells.Item(19, 1).Value = 'knee pad cover'
$ws.Cells.Item(20, 1).Value = 'knee pad cycling'
$ws.Cells.Item(21, 1).Value = 'knee pad inserts for work pants'
$ws.Cells.Item(22, 1).Value = 'knee pad pants women'
$ws.Cells.Item(23, 1).Value = 'knee pad pink'
$ws.Cells.Item(24, 1).Value = 'knee pad scooter'
$ws.Cells.Item(25, 1).Value = 'knee pad skating'
$ws.Cells.Item(26, 1).Value = 'knee pad wheels'
$ws.Cells.Item(27, 1).Value = 'neoprene knee pads'
$ws.Cells.Item(28, 1).Value = 'skating knee pads'
$ws.Cells.Item(29, 1).Value = 'knee pads capezio'
$ws.Cells.Item(30, 1).Value = 'knee pads canoe'
$ws.Cells.Item(31, 1).Value = 'knee pads cycle'
$ws.Cells.Item(32, 1).Value = 'knee pads firefighter'
$ws.Cells.Item(33, 1).Value = 'knee pads gymnastics'
$ws.Cells.Item(34, 1).Value = 'knee pads multicam'
$ws.Cells.Item(35, 1).Value = 'knee pads protec'
$ws.Cells.Item(36, 1).Value = 'knee pads shin'
$ws.Cells.Item(37, 1).Value = 'knee pads silver'
$ws.Cells.Item(38, 1).Value = 'skate knee pads'
$ws.Cells.Item(39, 1).Value = 'ski knee pad'
$ws.Cells.Item(40, 1).Value = 'kids knee protector'
$ws.Cells.Item(41, 1).Value = 'mens warming compression pants'
$ws.Cells.Item(42, 1).Value = 'mens workout tights'
$ws.Cells.Item(43, 1).Value = 'womens knee pads basketball'
$ws.Cells.Item(44, 1).Value = 'nike pro compression tights men'
$ws.Cells.Item(45, 1).Value = 'mens basketball pants'
$ws.Cells.Item(46, 1).Value = 'basketball knee pads for kids boys'
$ws.Cells.Item(47, 1).Value = 'nike basketball tights'
$ws.Cells.Item(48, 1).Value = 'military pants with knee pads'
$ws.Cells.Item(49, 1).Value = 'asics knee pads'
$ws.Cells.Item(50, 1).Value = 'knee pads xlarge'
$ws.Cells.Item(51, 1).Value = 'motorcycle knee pads men'
$ws.Cells.Item(52, 1).Value = 'gray baseball pants mens'
$ws.Cells.Item(53, 1).Value = 'baseball pants mens knickers'
$ws.Cells.Item(54, 1).Value = 'kids basketball knee pads youth'
$ws.Cells.Item(55, 1).Value = 'youth knee pads basketball for kids'
$ws.Cells.Item(56, 1).Value = 'workout leggings for men'
$ws.Cells.Item(57, 1).Value = 'legging for men nike'
$ws.Cells.Item(58, 1).Value = 'adidas capris men'
$ws.Cells.Item(59, 1).Value = 'youth knee and elbow pads'
$ws.Cells.Item(60, 1).Value = 'knee pads for dancers'
$ws.Cells.Item(61, 1).Value = 'knee pad bathtub'
$ws.Cells.Item(62, 1).Value = 'under armour compression tights men'
$ws.Cells.Item(63, 1).Value = 'knee pad for dancers'
$ws.Cells.Item(64, 1).Value = 'nike youth basketball tights'
$ws.Cells.Item(65, 1).Value = 'compression pants women'
$ws.Cells.Item(66, 1).Value = 'mens workout tights pants'
$ws.Cells.Item(67, 1).Value = 'white nike compression pants men'
$ws.Cells.Item(68, 1).Value = 'defender mens compression pants'
$ws.Cells.Item(69, 1).Value = 'mens nike basketball pants'
$ws.Cells.Item(70, 1).Value = 'nike basketball pants men'
$ws.Cells.Item(71, 1).Value = 'mens leggings compression nike'
$ws.Cells.Item(72, 1).Value = 'white compression pants men'
$ws.Cells.Item(73, 1).Value = 'od green pants with knee pads'
$ws.Cells.Item(74, 1).Value = 'send knee pad'
$ws.Cells.Item(75, 1).Value = 'children knee pads'
$ws.Cells.Item(76, 1).Value = 'smith knee pads'
$ws.Cells.Item(77, 1).Value = 'youth xl football pants'
$ws.Cells.Item(78, 1).Value = 'labor knee pads'
$ws.Cells.Item(79, 1).Value = 'rubber knee pads'
$ws.Cells.Item(80, 1).Value = 'compression pants nike'
$ws.Cells.Item(81, 1).Value = 'athletic capris for women'
$ws.Cells.Item(82, 1).Value = 'skins tights men'
$ws.Cells.Item(83, 1).Value = 'supportive knee pads'
$ws.Cells.Item(84, 1).Value = 'dye knee pads'
$ws.Cells.Item(85, 1).Value = 'dancers knee pads'
$ws.Cells.Item(86, 1).Value = 'apex knee pads'
$ws.Cells.Item(87, 1).Value = 'elbow and knee pads'
$ws.Cells.Item(88, 1).Value = 'fuse knee pads'
$ws.Cells.Item(89, 1).Value = 'ama knee pads'
$ws.Cells.Item(90, 1).Value = 'kp knee pads'
$ws.Cells.Item(91, 1).Value = 'adidas tights men'
$ws.Cells.Item(92, 1).Value = 'caterpillar knee pads'
$ws.Cells.Item(93, 1).Value = 'husky knee pads'
$ws.Cells.Item(94, 1).Value = 'mens leggings white'
$ws.Cells.Item(95, 1).Value = 'nike youth compression pants'
$ws.Cells.Item(96, 1).Value = 'muscle leggings men'
$ws.Cells.Item(97, 1).Value = 'mens compression pants under armour'
$ws.Cells.Item(98, 1).Value = 'airsoft pants with knee pads'
$ws.Cells.Item(99, 1).Value = 'football pants youth with pads'
$ws.Cells.Item(100, 1).Value = 'lotus leggings men'
